# Weekly refresh of the "Chirimoya" price sheet: a new, more recent
# weekly record is inserted at the top of the data block (row 4), pushing
# the existing rows 4-13 down to rows 5-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 (data for row 4 onward shifts
# down by one; row 3's formatting - e.g. the date style on column D - is
# inherited by the newly inserted row).
$ws.Rows(4).Insert()

# Populate the new row 4 with this week's record.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44811
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107002
$ws.Range("J4").Value = "Chirimoya"
$ws.Range("K4").Value = "Cultivar IV Región"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 29000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 29500
$ws.Range("Q4").Value = "`$/caja 12 kilos"
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 2458
$ws.Range("T4").Value = 12
